# Burndown - Sprint 4: Meeting Recording & Burndown
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New task row "Other Tasks" added at row 23 with 5 estimated hours,
# and 5 hours actually burned down on Day 2 (column E).
$ws.Range("A23").Value = "Other Tasks"
$ws.Range("B23").Value = 5
$ws.Range("E23").Value = 5

# Existing task on row 12 (Input Manager Sortout) also had 2 hours burned
# down on Day 2 (column E).
$ws.Range("E12").Value = 2

# Update the active cell selection to reflect where the user was working.
$ws.Range("E15").Select()

$wb.Save()
